$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MT"
$ws.Range("D2").Value = 285.635
$ws.Range("E2").Value = 1274.1262
$ws.Range("F2").Value = 354.76
$ws.Range("G2").Value = 0.2241810897539035
$ws.Range("H2").Value = 209321758.9299
$ws.Range("I2").Value = -1.9552
